$wb = $excel.ActiveWorkbook

# --- Sheet1: reserva_total.prn ---
$ws1 = $wb.Worksheets.Item("reserva_total.prn")
$ws1.Range("B1:D1").ClearContents()
$ws1.Range("A1").Value = "Análisis de la Reserva Total"

$ws1.Range("A3").Value = "RESERVA ROTANTE EN MAQUINAS QUE REGULAN"
$ws1.Range("A3:F3").Merge()
$ws1.Range("A3:F3").HorizontalAlignment = -4108
$ws1.Range("A3:F3").Borders.LineStyle = 1

$ws1.Range("A4").Value = "RESERVA HIDRO [MW]"
$ws1.Range("A4:C4").Merge()
$ws1.Range("D4:F4").Merge()

$ws1.Range("A5").Value = "RESERVA TERMICA [MW]"
$ws1.Range("A5:C5").Merge()
$ws1.Range("D5:F5").Merge()

$ws1.Range("A6").Value = "RESERVA TOTAL [MW]"
$ws1.Range("A6:C6").Merge()
$ws1.Range("D6:F6").Merge()

$ws1.Range("A7").Value = "RESERVA ROTANTE DEL PARQUE REGULANTE"
$ws1.Range("A7:F7").Merge()
$ws1.Range("A7:F7").HorizontalAlignment = -4108
$ws1.Range("A7:F7").Borders.LineStyle = 1

$ws1.Range("A8").Value = "RESERVA HIDRO"

$ws1.Range("A9").Value = "RESERVA PROGRAMADA A 50Hz PARA RPF"

$ws1.Range("A10").Value = "RESERVA HIDRO"

$ws1.Range("A11").Value = "RESERVA TÉRMICA"

$ws1.Range("A12").Value = "TOTAL SISTEMA"

$ws1.Range("A13").Value = "RESERVA PARA RPF"

$ws1.Range("A14").Value = "COLABORACIÓN DEL PARQUE HIDRO EN RSF [MW]"

$ws1.Range("A15").Value = "COLABORACIÓN DEL PARQUE HIDRO EN RSF [%]"

$ws1.Range("A16").Value = "POTENCIA OPERABLE EN EL PARQUE REGULANTE"
$ws1.Range("A16:F16").Merge()
$ws1.Range("A16:F16").HorizontalAlignment = -4108
$ws1.Range("A16:F16").Borders.LineStyle = 1

$ws1.Range("A17").Value = "HIDRO"

$ws1.Range("A18").Value = "TÉRMICA TG-CC"

$ws1.Range("A19").Value = "TÉRMICA TV"

$ws1.Range("A20").Value = "TOTAL"

$ws1.Range("A21").Value = "RESERVA PROGRAMADA EN EL PARQUE REGULANTE"
$ws1.Range("A21:F21").Merge()
$ws1.Range("A21:F21").HorizontalAlignment = -4108
$ws1.Range("A21:F21").Borders.LineStyle = 1

$ws1.Range("A22").Value = "HIDRO"

$ws1.Range("A23").Value = "TÉRMICA TG-CC"

$ws1.Range("A24").Value = "TÉRMICA TV"

$ws1.Range("A25").Value = "TOTAL"

$ws1.Range("A26").Value = "RESERVA NUEVA"

$ws1.Range("A27").Value = "RESERVA TOTAL 2"

# --- Sheet6: Reserva.err ---
$ws6 = $wb.Worksheets.Item("Reserva.err")
$ws6.Range("A2").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98011 SGDEHI0713.8"
$ws6.Range("A3").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98011 SGDEHI0713.8"
$ws6.Range("A4").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98012 SGDEHI0813.8"
$ws6.Range("A5").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98012 SGDEHI0813.8"
$ws6.Range("A6").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98013 SGDEHI0913.8"
$ws6.Range("A7").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98013 SGDEHI0913.8"
$ws6.Range("A8").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98014 SGDEHI1013.8"
$ws6.Range("A9").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98014 SGDEHI1013.8"
$ws6.Range("A10").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98015 SGDEHI1113.8"
$ws6.Range("A11").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98015 SGDEHI1113.8"
$ws6.Range("A12").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98016 SGDEHI1213.8"
$ws6.Range("A13").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98016 SGDEHI1213.8"
$ws6.Range("A14").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98017  SGDEHI1413.8"
$ws6.Range("A15").Value = "***** ERROR EN LOS DATOS DE GENSALE.PRN ***** NO SE ENCUENTRA LA BARRA 98017  SGDEHI1413.8"
$ws6.Range("A16").Value = "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 99 U.T.E NO POSEE SYSTEMA"
$ws6.Range("A17").Value = "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 18 A.N.D.E. NO POSEE SYSTEMA"
$ws6.Range("A18").Value = "***** ERROR EN LOS DATOS DE reserva_DEMANDAS ***** EL AREA INDICADA COMO 20 BRASIL NO POSEE SYSTEMA"
